# Eliminada columna de valores calibrados para Temp y rellenado de tres
# celdas de temperatura (sheet "MEDIDAS_ADC_A10_CALIBRAR(EC12b2")
#
# Summary of the change:
#  - The hidden "Nadc_corr" helper column (J, containing
#    H*(C10/2^15)*(C.. /2^15)+C9 formulas) is removed entirely; the
#    "Tcalib" column that was in K shifts left into J, and the "Mean"
#    column that was in L shifts left into K (its SUM formula now points
#    at J instead of K). The results table that was in column N (labels +
#    the "Medidas tras la calibración" header) shifts left into M.
#  - Three previously-empty cells in the results column (M9, M11, M13)
#    are filled in with the corresponding calibrated-mean values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MEDIDAS_ADC_A10_CALIBRAR(EC12b2")

# Remove the whole "Nadc_corr" column (J). Everything to its right
# (Tcalib/Mean in K:L, and the results block in N) shifts one column to
# the left automatically, and all formulas/shared-strings are re-wired by
# Excel as part of the delete.
$ws.Columns.Item(10).EntireColumn.Delete() | Out-Null

# Fill in the three blank "Medida ADC ºC" result cells with the mean
# calibrated-temperature values computed a few rows below/above them.
$ws.Range("M9").Formula = "=K9"
$ws.Range("M11").Formula = "=K14"
$ws.Range("M13").Formula = "=K19"

# The now-visible column I (former hidden helper column position) is
# widened.
$ws.Columns.Item(9).ColumnWidth = 17.333333333333332

# Restore the user's selection/active cell.
$ws.Range("M14").Select() | Out-Null
